$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's name columns (C/D) use a trailing "<space><nbsp>" pair
# (matches the source file's existing convention, e.g. "Marni &#160;").
$nbsp = [char]0x00A0

# These cells hold purely numeric-looking text (e.g. "2", "22.0") that must
# stay stored as TEXT (matching the source file's t="inlineStr" convention)
# rather than being auto-coerced to a number by Excel. Forcing a Text
# number format keeps the literal string (incl. trailing ".0") intact.
$numericLookingTextCells = @(
    "B6","H6",
    "B7","H7",
    "B8","H8",
    "B9","H9",
    "B10","H10",
    "B13"
)
foreach ($addr in $numericLookingTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 6
$ws.Range("B6").Value = "2"
$ws.Range("C6").Value = "Elwanda " + $nbsp
$ws.Range("D6").Value = "Cassy " + $nbsp
$ws.Range("E6").Value = "-1.07,-9.07"
$ws.Range("F6").Value = "Tamisha(mother): 0550693864"
$ws.Range("H6").Value = "22.0"

# Row 7
$ws.Range("B7").Value = "0"
$ws.Range("C7").Value = "Trudie " + $nbsp
$ws.Range("D7").Value = "Fleta " + $nbsp
$ws.Range("E7").Value = "-0.75,-7.75"
$ws.Range("F7").Value = "Anneliese(father): 0548973345"
$ws.Range("G7").Value = "7:02:00"
$ws.Range("H7").Value = "20.0"

# Row 8
$ws.Range("B8").Value = "8"
$ws.Range("C8").Value = "Marni " + $nbsp
$ws.Range("D8").Value = "Shanika " + $nbsp
$ws.Range("E8").Value = "1.1,-7.16"
$ws.Range("F8").Value = "Lady(mother): 0560804012"
$ws.Range("G8").Value = "7:05:00"
$ws.Range("H8").Value = "17.0"

# Row 9
$ws.Range("B9").Value = "12"
$ws.Range("C9").Value = "Frankie " + $nbsp
$ws.Range("D9").Value = "Flavia " + $nbsp
$ws.Range("E9").Value = "-2.25,-1.67"
$ws.Range("F9").Value = "Cyrus(mother): 0522363358"
$ws.Range("G9").Value = "7:14:00"
$ws.Range("H9").Value = "8.0"

# Row 10
$ws.Range("B10").Value = "11"
$ws.Range("C10").Value = "Randolph " + $nbsp
$ws.Range("D10").Value = "Bridgette " + $nbsp
$ws.Range("E10").Value = "1.09,-0.75"
$ws.Range("F10").Value = "Lenny(father): 0505536740"
$ws.Range("G10").Value = "7:20:00"
$ws.Range("H10").Value = "2.0"

# Row 11 (school)
$ws.Range("G11").Value = "7:22:00"

# Row 13 (time)
$ws.Range("B13").Value = "22.0"
